$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "paulo mandou abaixar o volume"
$ws.Range("B10").Value = "06/09/2025 23:47:00"

$ws.Range("A11").Value = "testandooo"
$ws.Range("B11").Value = "06/09/2025 23:47:43"

$ws.Range("A12").Value = "testeee oii"
$ws.Range("B12").Value = "06/09/2025 23:47:59"
